$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Plxna1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.339548666666667
$ws.Range("H2").Value = 4.018646
$ws.Range("I2").Value = 0.09827984122213275
$ws.Range("J2").Value = 0.09827984122213274
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.145781666666667
$ws.Range("N2").Value = 21.437345
$ws.Range("O2").Value = 0.1148763047483796
$ws.Range("P2").Value = 0.1148763047483796
$ws.Range("Q2").Value = 9.572122303874446
$ws.Range("R2").Value = 86.14910073487
$ws.Range("S2").Value = 0.01129002499085608
$ws.Range("T2").Value = 0.01129002499085608

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Plxna1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.339548666666667
$ws.Range("H3").Value = 4.018646
$ws.Range("I3").Value = 0.09827984122213275
$ws.Range("J3").Value = 0.09827984122213274
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.53178066666667
$ws.Range("N3").Value = 55.595342
$ws.Range("O3").Value = 0.2979187698001963
$ws.Range("P3").Value = 0.2979187698001963
$ws.Range("Q3").Value = 24.82422208299245
$ws.Range("R3").Value = 223.417998746932
$ws.Range("S3").Value = 0.02927940939305641
$ws.Range("T3").Value = 0.0292794093930564

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Plxna1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.339548666666667
$ws.Range("H4").Value = 4.018646
$ws.Range("I4").Value = 0.09827984122213275
$ws.Range("J4").Value = 0.09827984122213274
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.83059133333333
$ws.Range("N4").Value = 68.49177399999999
$ws.Range("O4").Value = 0.3670268824232265
$ws.Range("P4").Value = 0.3670268824232265
$ws.Range("Q4").Value = 30.58268817977822
$ws.Range("R4").Value = 275.244193618004
$ws.Range("S4").Value = 0.03607134372880909
$ws.Range("T4").Value = 0.03607134372880908

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Plxna1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.339548666666667
$ws.Range("H5").Value = 4.018646
$ws.Range("I5").Value = 0.09827984122213275
$ws.Range("J5").Value = 0.09827984122213274
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.69598566666667
$ws.Range("N5").Value = 41.087957
$ws.Range("O5").Value = 0.2201780430281976
$ws.Range("P5").Value = 0.2201780430281976
$ws.Range("Q5").Value = 18.34643933846911
$ws.Range("R5").Value = 165.117954046222
$ws.Range("S5").Value = 0.02163906310941117
$ws.Range("T5").Value = 0.02163906310941117

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Plxna1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.040291
$ws.Range("H6").Value = 6.120873
$ws.Range("I6").Value = 0.1496918182345096
$ws.Range("J6").Value = 0.1496918182345096
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.145781666666667
$ws.Range("N6").Value = 21.437345
$ws.Range("O6").Value = 0.1148763047483796
$ws.Range("P6").Value = 0.1148763047483796
$ws.Range("Q6").Value = 14.579474022465
$ws.Range("R6").Value = 131.215266202185
$ws.Range("S6").Value = 0.01719604292984658
$ws.Range("T6").Value = 0.01719604292984658

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Plxna1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.040291
$ws.Range("H7").Value = 6.120873
$ws.Range("I7").Value = 0.1496918182345096
$ws.Range("J7").Value = 0.1496918182345096
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.53178066666667
$ws.Range("N7").Value = 55.595342
$ws.Range("O7").Value = 0.2979187698001963
$ws.Range("P7").Value = 0.2979187698001963
$ws.Range("Q7").Value = 37.810225308174
$ws.Range("R7").Value = 340.292027773566
$ws.Range("S7").Value = 0.0445960023375797
$ws.Range("T7").Value = 0.0445960023375797

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Plxna1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.040291
$ws.Range("H8").Value = 6.120873
$ws.Range("I8").Value = 0.1496918182345096
$ws.Range("J8").Value = 0.1496918182345096
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.83059133333333
$ws.Range("N8").Value = 68.49177399999999
$ws.Range("O8").Value = 0.3670268824232265
$ws.Range("P8").Value = 0.3670268824232265
$ws.Range("Q8").Value = 46.58105002207799
$ws.Range("R8").Value = 419.2294501987019
$ws.Range("S8").Value = 0.05494092137087636
$ws.Range("T8").Value = 0.05494092137087636

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Plxna1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.040291
$ws.Range("H9").Value = 6.120873
$ws.Range("I9").Value = 0.1496918182345096
$ws.Range("J9").Value = 0.1496918182345096
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 13.69598566666667
$ws.Range("N9").Value = 41.087957
$ws.Range("O9").Value = 0.2201780430281976
$ws.Range("P9").Value = 0.2201780430281976
$ws.Range("Q9").Value = 27.943796291829
$ws.Range("R9").Value = 251.494166626461
$ws.Range("S9").Value = 0.032958851596207
$ws.Range("T9").Value = 0.032958851596207

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Plxna1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.25010366666667
$ws.Range("H10").Value = 30.750311
$ws.Range("I10").Value = 0.7520283405433575
$ws.Range("J10").Value = 0.7520283405433575
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.145781666666667
$ws.Range("N10").Value = 21.437345
$ws.Range("O10").Value = 0.1148763047483796
$ws.Range("P10").Value = 0.1148763047483796
$ws.Range("Q10").Value = 73.24500286269945
$ws.Range("R10").Value = 659.205025764295
$ws.Range("S10").Value = 0.08639023682767695
$ws.Range("T10").Value = 0.08639023682767695

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema3a"
$ws.Range("C11").Value = "Plxna1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.25010366666667
$ws.Range("H11").Value = 30.750311
$ws.Range("I11").Value = 0.7520283405433575
$ws.Range("J11").Value = 0.7520283405433575
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 18.53178066666667
$ws.Range("N11").Value = 55.595342
$ws.Range("O11").Value = 0.2979187698001963
$ws.Range("P11").Value = 0.2979187698001963
$ws.Range("Q11").Value = 189.9526729612624
$ws.Range("R11").Value = 1709.574056651362
$ws.Range("S11").Value = 0.2240433580695601
$ws.Range("T11").Value = 0.2240433580695601

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema3a"
$ws.Range("C12").Value = "Plxna1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.25010366666667
$ws.Range("H12").Value = 30.750311
$ws.Range("I12").Value = 0.7520283405433575
$ws.Range("J12").Value = 0.7520283405433575
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 22.83059133333333
$ws.Range("N12").Value = 68.49177399999999
$ws.Range("O12").Value = 0.3670268824232265
$ws.Range("P12").Value = 0.3670268824232265
$ws.Range("Q12").Value = 234.0159279379682
$ws.Range("R12").Value = 2106.143351441714
$ws.Range("S12").Value = 0.276014617323541
$ws.Range("T12").Value = 0.276014617323541

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema3a"
$ws.Range("C13").Value = "Plxna1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.25010366666667
$ws.Range("H13").Value = 30.750311
$ws.Range("I13").Value = 0.7520283405433575
$ws.Range("J13").Value = 0.7520283405433575
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 13.69598566666667
$ws.Range("N13").Value = 41.087957
$ws.Range("O13").Value = 0.2201780430281976
$ws.Range("P13").Value = 0.2201780430281976
$ws.Range("Q13").Value = 140.3852729005141
$ws.Range("R13").Value = 1263.467456104627
$ws.Range("S13").Value = 0.1655801283225794
$ws.Range("T13").Value = 0.1655801283225794

$ws.Rows("14:17").Delete()
